$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (this is what drives the Print_Area defined name update below)
$ws.Name = "Uncut_Sheet_1"

# Re-apply the print area so the defined name reflects the new sheet name while
# keeping the same printed range ($A$1:$G$42)
$ws.PageSetup.PrintArea = "`$A`$1:`$G`$42"

# Update the active selection on the sheet to the merged range B16:C16
$ws.Range("B16:C16").Select()
